$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $value)
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

# Row 2
Set-TextValue "D2" "260.93"
Set-TextValue "E2" "0.01%"
# Row 3
Set-TextValue "D3" "26.86"
# Row 4
Set-TextValue "D4" "4.715"
Set-TextValue "E4" "0.10%"
# Row 5
Set-TextValue "D5" "0.06223"
Set-TextValue "E5" "2.37%"
# Row 6
Set-TextValue "D6" "6.732"
Set-TextValue "E6" "0.89%"
# Row 7
Set-TextValue "D7" "0.8493"
Set-TextValue "E7" "0.38%"
# Row 8
Set-TextValue "D8" "0.9112"
Set-TextValue "E8" "-1.17%"
# Row 9
Set-TextValue "D9" "0.1401"
Set-TextValue "E9" "-0.14%"
# Row 10
Set-TextValue "D10" "0.04944"
Set-TextValue "E10" "-1.60%"
# Row 11
Set-TextValue "D11" "0.07073"
Set-TextValue "E11" "-0.35%"
# Row 12
Set-TextValue "D12" "0.03085"
Set-TextValue "E12" "-1.41%"
# Row 13
Set-TextValue "D13" "0.09055"
Set-TextValue "E13" "-0.19%"
# Row 14
Set-TextValue "D14" "0.001530"
Set-TextValue "E14" "-0.44%"
# Row 15
Set-TextValue "D15" "0.0006187"
Set-TextValue "E15" "1.39%"
# Row 16
Set-TextValue "D16" "0.005963"
Set-TextValue "E16" "-2.79%"
# Row 17
Set-TextValue "E17" "-0.15%"
# Row 18
Set-TextValue "E18" "0.84%"
# Row 19
Set-TextValue "D19" "2.167"
Set-TextValue "E19" "0.02%"
# Row 21
Set-TextValue "E21" "1.04%"
# Row 22
Set-TextValue "D22" "4.089"
Set-TextValue "E22" "-0.08%"
# Row 23
Set-TextValue "D23" "0.04247"
Set-TextValue "E23" "0.29%"
# Row 24
Set-TextValue "D24" "0.001202"
Set-TextValue "E24" "-1.57%"
# Row 25
Set-TextValue "D25" "0.004075"
Set-TextValue "E25" "4.18%"
# Row 26
Set-TextValue "E26" "-0.02%"
# Row 27
Set-TextValue "E27" "4.07%"
# Row 40
Set-TextValue "D40" "0.03940"
Set-TextValue "E40" "1.75%"
# Row 41
Set-TextValue "E41" "0.00%"
# Row 42
Set-TextValue "D42" "0.004132"
# Row 43
Set-TextValue "E43" "-6.13%"
# Row 44
Set-TextValue "D44" "0.01339"
Set-TextValue "E44" "-18.13%"
# Row 45
Set-TextValue "E45" "-2.97%"
# Row 46
Set-TextValue "E46" "-0.02%"
# Row 47
Set-TextValue "E47" "-37.60%"
# Row 48
Set-TextValue "D48" "0.2524"
Set-TextValue "E48" "85.83%"
# Row 49
Set-TextValue "E49" "-0.02%"
# Row 50
Set-TextValue "E50" "-0.02%"
